$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Description sheet: wording/number refresh + new footnote row
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Description")
$ws2.Activate()

# Update the "Desc Body" paragraph: 2015 -> 2016, 3,272 -> 3,269
$ws2.Range("B5").Value = "From 2008 to the end of October 2016, the NPRH and NPARIH delivered a total of 3,269 new houses against the 2018 COAG target of 4,200."
$ws2.Rows.Item(5).RowHeight = 26.85

# Insert a new footnote row right before the ACT/Vic/Source block and
# give it the long-text note wording + matching row height.
$ws2.Rows.Item(10).Insert()
$ws2.Range("B10").Value = "Totals include dwellings built from Commonwealth Own Purpose Expense funding under NPARIH."
$ws2.Rows.Item(10).RowHeight = 26.95

# Park the view where the diff says it ended up (before Data becomes active).
$ws2.Range("B10").Select() | Out-Null

# ---------------------------------------------------------------------
# Data sheet: refreshed 2016 totals
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Data")
$ws1.Activate()
$ws1.Range("D4").Value = 833
$ws1.Range("I4").Value = 3269
$ws1.Range("I5").Select() | Out-Null
